# Generate Report for Handback
# - Status "Ready for handoff" -> "Handback transform failed" (Overview, zh-cn, de-de)
# - Error Detail column (P) widened to fit the new handback-mismatch message
# - Error Detail populated for the 348c437f... row on both the zh-cn and de-de sheets

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status changed from "Ready for handoff" to "Handback transform failed"
# for the 348c437f-23c7-4024-9874-a2a06b50c656 row, everywhere it's shown.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Error Detail (column P) widened from ~13.75 to 40 characters on both
# language sheets. ColumnWidth is offset from the stored sheet-XML width by
# 5/6 of a character, so subtract that to land exactly on 40.
$zhcn.Columns.Item(16).ColumnWidth = 40 - 5/6
$dede.Columns.Item(16).ColumnWidth = 40 - 5/6

# Error Detail message for the failed handback transform, one per locale.
$zhcn.Range("P3").Value = "Handback file name: bes2bysz.dyi is different with handoff file name: 348c437f-23c7-4024-9874-a2a06b50c656.4435ef1714ce749182b12372dd50a468ccde892e.zh-cn."
$dede.Range("P3").Value = "Handback file name: bes2bysz.dyi is different with handoff file name: 348c437f-23c7-4024-9874-a2a06b50c656.4435ef1714ce749182b12372dd50a468ccde892e.de-de."
